$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026037100109277
$ws.Range("D2").Value = 1.03005174074101
$ws.Range("E2").Value = 1.026320117792564
$ws.Range("I2").Value = 1.032549613526556
$ws.Range("J2").Value = 1.031202777825712
$ws.Range("K2").Value = 1.032863871900083
$ws.Range("L2").Value = 1.02914310931884
$ws.Range("N2").Value = 1.032667204050842
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.0267905934845
$ws.Range("D3").Value = 1.030587159621844
$ws.Range("E3").Value = 1.026954644493895
$ws.Range("I3").Value = 1.032676756920272
$ws.Range("J3").Value = 1.031596870503527
$ws.Range("K3").Value = 1.033208298636306
$ws.Range("L3").Value = 1.029585594238767
$ws.Range("N3").Value = 1.033061856385457
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027278776838372
$ws.Range("D4").Value = 1.030934093976628
$ws.Range("E4").Value = 1.027366159003453
$ws.Range("I4").Value = 1.03275811027529
$ws.Range("J4").Value = 1.031851812767169
$ws.Range("K4").Value = 1.033430937701247
$ws.Range("L4").Value = 1.029872155688618
$ws.Range("N4").Value = 1.033317160696356
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027484156211686
$ws.Range("D5").Value = 1.031080058759029
$ws.Range("E5").Value = 1.027539381090952
$ws.Range("I5").Value = 1.032792090962818
$ws.Range("J5").Value = 1.031958974245698
$ws.Range("K5").Value = 1.033524479218846
$ws.Range("L5").Value = 1.029992683042867
$ws.Range("N5").Value = 1.033424474356476
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027518648879618
$ws.Range("D6").Value = 1.031104573492119
$ws.Range("E6").Value = 1.027568478764328
$ws.Range("I6").Value = 1.032797783531667
$ws.Range("J6").Value = 1.031976966136061
$ws.Range("K6").Value = 1.033540181931418
$ws.Range("L6").Value = 1.030012923409556
$ws.Range("N6").Value = 1.033442491797387
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027281520552146
$ws.Range("D7").Value = 1.030936043921352
$ws.Range("E7").Value = 1.027368472738143
$ws.Range("I7").Value = 1.032758565193318
$ws.Range("J7").Value = 1.03185324472969
$ws.Range("K7").Value = 1.033432187829103
$ws.Range("L7").Value = 1.02987376595928
$ws.Range("N7").Value = 1.033318594692428
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026291616576448
$ws.Range("D8").Value = 1.030232586980806
$ws.Range("E8").Value = 1.026534364585581
$ws.Range("I8").Value = 1.032592771635563
$ws.Range("J8").Value = 1.031335975301738
$ws.Range("K8").Value = 1.032980318941569
$ws.Range("L8").Value = 1.02929259770862
$ws.Range("N8").Value = 1.032800590682556
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024552137262274
$ws.Range("D9").Value = 1.028996787340949
$ws.Range("E9").Value = 1.025071797908159
$ws.Range("I9").Value = 1.032293632875526
$ws.Range("J9").Value = 1.030424066865344
$ws.Range("K9").Value = 1.03218238102849
$ws.Range("L9").Value = 1.028270441906977
$ws.Range("N9").Value = 1.031887387231586
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023395870090867
$ws.Range("D10").Value = 1.028175579302719
$ws.Range("E10").Value = 1.024101742639391
$ws.Range("I10").Value = 1.032089552736406
$ws.Range("J10").Value = 1.029815933159764
$ws.Range("K10").Value = 1.031649368538103
$ws.Range("L10").Value = 1.027590397965151
$ws.Range("N10").Value = 1.03127838990639
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022896021981507
$ws.Range("D11").Value = 1.027820642316063
$ws.Range("E11").Value = 1.023682907367542
$ws.Range("I11").Value = 1.032000089360907
$ws.Range("J11").Value = 1.029552576022875
$ws.Range("K11").Value = 1.031418335343552
$ws.Range("L11").Value = 1.027296281192752
$ws.Range("N11").Value = 1.031014658772159
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022710481867348
$ws.Range("D12").Value = 1.027688902807964
$ws.Range("E12").Value = 1.023527516472392
$ws.Range("I12").Value = 1.031966694878572
$ws.Range("J12").Value = 1.029454750196236
$ws.Range("K12").Value = 1.031332485364403
$ws.Range("L12").Value = 1.027187086792539
$ws.Range("N12").Value = 1.030916694021625
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022750275129198
$ws.Range("D13").Value = 1.027717156827102
$ws.Range("E13").Value = 1.023560840023811
$ws.Range("I13").Value = 1.031973865508549
$ws.Range("J13").Value = 1.029475734292682
$ws.Range("K13").Value = 1.031350901991121
$ws.Range("L13").Value = 1.027210506908432
$ws.Range("N13").Value = 1.030937707917895
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022880682601081
$ws.Range("D14").Value = 1.027809750637851
$ws.Range("E14").Value = 1.023670058958711
$ws.Range("I14").Value = 1.031997332299682
$ws.Range("J14").Value = 1.029544489767758
$ws.Range("K14").Value = 1.03141123964236
$ws.Range("L14").Value = 1.027287254045753
$ws.Range("N14").Value = 1.031006561033632
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022961047670247
$ws.Range("D15").Value = 1.027866814037356
$ws.Range("E15").Value = 1.023737376684966
$ws.Range("I15").Value = 1.032011769281665
$ws.Range("J15").Value = 1.029586851895979
$ws.Range("K15").Value = 1.031448411201447
$ws.Range("L15").Value = 1.027334547648038
$ws.Range("N15").Value = 1.031048983320934
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02342906088283
$ws.Range("D16").Value = 1.028199149199883
$ws.Range("E16").Value = 1.024129564927879
$ws.Range("I16").Value = 1.032095467113491
$ws.Range("J16").Value = 1.029833410769892
$ws.Range("K16").Value = 1.031664696611555
$ws.Range("L16").Value = 1.027609924966269
$ws.Range("N16").Value = 1.031295892336729
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023722855050418
$ws.Range("D17").Value = 1.028407790116221
$ws.Range("E17").Value = 1.024375898383047
$ws.Range("I17").Value = 1.032147675837439
$ws.Range("J17").Value = 1.029988063287888
$ws.Range("K17").Value = 1.031800304723636
$ws.Range("L17").Value = 1.027782755824914
$ws.Range("N17").Value = 1.031450764479039
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023894299507586
$ws.Range("D18").Value = 1.028529549559864
$ws.Range("E18").Value = 1.024519696648366
$ws.Range("I18").Value = 1.032178022568529
$ws.Range("J18").Value = 1.030078266376028
$ws.Range("K18").Value = 1.031879379848472
$ws.Range("L18").Value = 1.027883598491105
$ws.Range("N18").Value = 1.03154109566591
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023952771023465
$ws.Range("D19").Value = 1.028571076961975
$ws.Range("E19").Value = 1.024568747773403
$ws.Range("I19").Value = 1.032188352062138
$ws.Range("J19").Value = 1.030109022727081
$ws.Range("K19").Value = 1.031906338514437
$ws.Range("L19").Value = 1.027917988843253
$ws.Range("N19").Value = 1.03157189569451
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023691325500565
$ws.Range("D20").Value = 1.028385398417806
$ws.Range("E20").Value = 1.024349457107528
$ws.Range("I20").Value = 1.03214208526658
$ws.Range("J20").Value = 1.029971470847588
$ws.Range("K20").Value = 1.031785757604173
$ws.Range("L20").Value = 1.027764209248522
$ws.Range("N20").Value = 1.031434148475571
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022842277363705
$ws.Range("D21").Value = 1.027782481294795
$ws.Range("E21").Value = 1.023637891631055
$ws.Range("I21").Value = 1.031990426430204
$ws.Range("J21").Value = 1.029524243072537
$ws.Range("K21").Value = 1.031393472634426
$ws.Range("L21").Value = 1.027264652432681
$ws.Range("N21").Value = 1.030986285585781
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022309175548938
$ws.Range("D22").Value = 1.027403982306817
$ws.Range("E22").Value = 1.023191563008209
$ws.Range("I22").Value = 1.031894125172939
$ws.Range("J22").Value = 1.029243035065816
$ws.Range("K22").Value = 1.031146631682071
$ws.Range("L22").Value = 1.026950872497322
$ws.Range("N22").Value = 1.030704678231436
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022591713173873
$ws.Range("D23").Value = 1.027604576235582
$ws.Range("E23").Value = 1.023428068992368
$ws.Range("I23").Value = 1.031945265833438
$ws.Range("J23").Value = 1.029392109994007
$ws.Range("K23").Value = 1.031277504841903
$ws.Range("L23").Value = 1.027117183150271
$ws.Range("N23").Value = 1.030853964863124
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.02370557210225
$ws.Range("D24").Value = 1.028395516068939
$ws.Range("E24").Value = 1.024361404423259
$ws.Range("I24").Value = 1.032144611731711
$ws.Range("J24").Value = 1.029978968266946
$ws.Range("K24").Value = 1.031792330891897
$ws.Range("L24").Value = 1.027772589545232
$ws.Range("N24").Value = 1.031441656542124
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025001246356377
$ws.Range("D25").Value = 1.029315811086175
$ws.Range("E25").Value = 1.025449036023195
$ws.Range("I25").Value = 1.03237179106008
$ws.Range("J25").Value = 1.030659857476102
$ws.Range("K25").Value = 1.032388858745904
$ws.Range("L25").Value = 1.028534454763477
$ws.Range("N25").Value = 1.032123512692056
